$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the view: zoom to 200% and move the active selection to F15
$excel.ActiveWindow.Zoom = 200
$ws.Range("F15").Select()

# Add the numbered comment/legend entries below the data table (column A,
# rows 16-23). Rows 16-21 get the "short" row height (15pt, matching the
# data rows above); rows 22-23 keep the template's existing 15.75pt height.
$comments = @(
    "# 1 Crownpoint Healthcare Facility (CHCF)",
    "# 2 Thoreau Health Station (THS)",
    "# 3 Pueblo Pintado Health Center (PPHC)",
    "# 4 No Preference (NP)",
    "# 5 Administrative Time (Adm)",
    "# 6 Approved Leave (AL)",
    "# 7 Continuing Medical Education (CME)",
    "# 8 General Staff Meeting (GME)"
)

$row = 16
foreach ($comment in $comments) {
    $ws.Cells.Item($row, 1).Value = $comment
    if ($row -le 21) {
        $ws.Rows.Item($row).RowHeight = 15
    }
    $row = $row + 1
}
